$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add the new "P1 arm infill test" print-spec rows (5, 6, 7), one per infill
# pattern being tested, mirroring the layout of the existing data rows.
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "P1 arm infill test"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 0.2
$ws.Range("D5").Value = "3D cubic"
$ws.Range("E5").Value = 25
$ws.Range("F5").Value = 0.4
$ws.Range("G5").Value = 0.4
$ws.Range("H5").Value = 0.4

$ws.Range("A6").Value = "P1 arm infill test"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 0.2
$ws.Range("D6").Value = "3D honeycomb"
$ws.Range("E6").Value = 25
$ws.Range("F6").Value = 0.4
$ws.Range("G6").Value = 0.4
$ws.Range("H6").Value = 0.4

$ws.Range("A7").Value = "P1 arm infill test"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 0.2
$ws.Range("D7").Value = "triangles (3d?)"
$ws.Range("E7").Value = 25
$ws.Range("F7").Value = 0.4
$ws.Range("G7").Value = 0.4
$ws.Range("H7").Value = 0.4

# The numeric/spec columns (B..H) are centered, matching the rest of the table.
$ws.Range("B5:H7").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# Box the whole data block (rows 4-7) in thin borders, like the rest of the
# table (row 3 already had this treatment).
# ---------------------------------------------------------------------------
$ws.Range("B4:J7").Borders.ColorIndex = 1
$ws.Range("B4:J7").Borders.LineStyle = 1

$ws.Range("A4:A7").Borders.ColorIndex = 1
$ws.Range("A4:A7").Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# Row 8: an extra formatted (but still empty) row under Weight/Comments,
# ready for the next entry.
# ---------------------------------------------------------------------------
$ws.Range("I8").Borders.ColorIndex = 1
$ws.Range("I8").Borders.LineStyle = 1
$ws.Range("I8").Borders.LineStyle = -4142

$ws.Range("J8").Borders.ColorIndex = 1
$ws.Range("J8").Borders.LineStyle = 1
$ws.Range("J8").Borders.LineStyle = -4142

# ---------------------------------------------------------------------------
# Column A (Part) and D (Infill Type) need to widen to fit the new, longer
# text that was just typed in ("P1 arm infill test", "3D honeycomb", ...).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 13.385416666666666
$ws.Columns.Item(4).ColumnWidth = 12.498697916666666

# Leave the selection on the last cell touched.
[void]$ws.Range("I8").Select()
